# Re-run of output tables: update the "Top Marginal Income Rate" (column C)
# values on Sheet1 to the refreshed figures for the affected countries.
# Values are stored as literal text (e.g. "60.2%"), matching how the rest
# of the table's percentage-like columns are authored, so we force a Text
# number format before writing and then clear the format back off so the
# cell keeps the literal string without picking up Excel's automatic
# percentage-number conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "C6"  = "60.2%"
    "C9"  = "31.1%"
    "C11" = "32.4%"
    "C12" = "58.4%"
    "C13" = "55.6%"
    "C16" = "33.5%"
    "C17" = "44.4%"
    "C18" = "52.0%"
    "C20" = "52.8%"
    "C21" = "56.1%"
    "C22" = "47.4%"
    "C23" = "38.9%"
    "C24" = "24.0%"
    "C25" = "47.2%"
    "C27" = "52.3%"
    "C29" = "46.6%"
    "C30" = "39.9%"
    "C31" = "58.2%"
    "C32" = "35.0%"
    "C33" = "61.1%"
    "C35" = "60.1%"
    "C37" = "45.5%"
    "C38" = "47.0%"
    "C39" = "46.0%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
